# Demo spreadsheet: Added format page
# Adds two new worksheets ("S3" and "MATLABFormat") after the existing
# "Sheet1" / "Day 2" sheets, populates their content/styles, and makes
# "S3" the active sheet (mirrors the author's workbook after the edit).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the two new sheets, in order, at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$s3 = $wb.Worksheets.Add($null, $lastSheet)
$s3.Name = "S3"

$matlab = $wb.Worksheets.Add($null, $s3)
$matlab.Name = "MATLABFormat"

# ---------------------------------------------------------------------
# 2. "S3" sheet content
#    Write cells in the same order the original author must have typed
#    them in, so that newly-interned shared strings come out in the
#    same sequence.
# ---------------------------------------------------------------------

# Row 3: column labels
$s3.Range("D3").Value = "Label"
$s3.Range("E3").Value = "Label 2"
$s3.Range("F3").Value = "Label 3"
$s3.Range("G3").Value = "Label 4"

# Row 4: row label + first data row
$s3.Range("C4").Value = "Test"
$s3.Range("E4").Value = 1
$s3.Range("F4").Value = 10
$s3.Range("G4").Value = 100

# Title (row 2)
$s3.Range("F2").Value = "Title"

# Notes column (K)
$s3.Range("K2").Value = "1. Find column before missing"
$s3.Range("K3").Value = "2. ID row before missing"
$s3.Range("K4").Value = "3. When extend column, check adjacent cells"

# Remaining data rows 5-7
$s3.Range("E5").Value = 2
$s3.Range("F5").Value = 12
$s3.Range("G5").Value = 200

$s3.Range("E6").Value = 3
$s3.Range("F6").Value = 13
$s3.Range("G6").Value = 300

$s3.Range("E7").Value = 4
$s3.Range("F7").Value = 14
$s3.Range("G7").Value = 400

# Row 8: repeat of the "Test" label (re-uses the existing shared string)
$s3.Range("H8").Value = "Test"

# Column K is noticeably wider (holds the long notes above)
$s3.Columns.Item(11).ColumnWidth = 26.1666

# Select the same cell the author had selected
$null = $s3.Range("J14").Select()

# ---------------------------------------------------------------------
# 3. "MATLABFormat" sheet content
# ---------------------------------------------------------------------

# Header row (bold)
$matlab.Range("B1").Value = "Sheet"
$matlab.Range("A1").Value = "Name"
$matlab.Range("C1").Value = "TL"
$matlab.Range("D1").Value = "BR"
$matlab.Range("A1:D1").Font.Bold = $true

# Data row 2 (plain formatting)
$matlab.Range("A2").Value = "Demo1"
$matlab.Range("B2").Value = "S3"
$matlab.Range("C2").Value = "E3"
$matlab.Range("D2").Value = "G7"

# Data row 3 (explicit Arial 10 formatting) -- write C,D,B,A so new
# shared strings are interned in the same order as the target file.
$matlab.Range("C3").Value = "A4"
$matlab.Range("D3").Value = "C8"
$matlab.Range("B3").Value = "Day 2"
$matlab.Range("A3").Value = "Demo 2"
$matlab.Range("A3:D3").Font.Name = "Arial"
$matlab.Range("A3:D3").Font.Size = 10
$matlab.Range("A3:D3").Font.Color = 0

# Column A is a bit wider to fit the sheet names
$matlab.Columns.Item(1).ColumnWidth = 13.4974

# Selection the author had on this sheet
$null = $matlab.Range("A8").Select()

# ---------------------------------------------------------------------
# 4. Make "S3" the active sheet/tab (matches activeTab="2" in the diff)
# ---------------------------------------------------------------------
$s3.Activate()
